$wb = $excel.ActiveWorkbook

# Sheet "展览" (Exhibition) - sheet1.xml
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F9").Value = 9029
$ws1.Range("F10").Value = 826
$ws1.Range("F13").Value = 1029
$ws1.Range("F21").Value = 1151

# Sheet "全部类型" (All Types) - sheet4.xml
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F11").Value = 9029
$ws4.Range("F12").Value = 826
$ws4.Range("F15").Value = 1029
$ws4.Range("F23").Value = 1151
